$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 76.23077000000001
$ws.Range("I9").Value = 80
$ws.Range("K9").Value = 80
$ws.Range("M9").Value = 89
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H28").Value = 957.1429000000001
$ws.Range("I28").Value = 760
$ws.Range("K28").Value = 760
$ws.Range("M28").Value = -275
$ws.Range("H58").Value = 599.8333
$ws.Range("J58").Value = 716.8
$ws.Range("L58").Value = 2150.4
$ws.Range("N58").Value = -2450.4
$ws.Range("H94").Value = 4349
$ws.Range("I94").Value = 4018.8
$ws.Range("K94").Value = 4018.8
$ws.Range("M94").Value = -3567.8
$ws.Range("H113").Value = 2009.3
$ws.Range("J113").Value = 2198
$ws.Range("L113").Value = 2198
$ws.Range("N113").Value = -8706
$ws.Range("H116").Value = 6452.5
$ws.Range("H132").Value = 7984.857
$ws.Range("I132").Value = 7984.857
$ws.Range("K132").Value = 23954.571
$ws.Range("M132").Value = -21424.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1250
$ws.Range("I35").Value = 1250
$ws.Range("K35").Value = 1250
$ws.Range("M35").Value = -844
$ws.Range("H97").Value = 736.6667
$ws.Range("I97").Value = 736.6667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 736.6667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -240.6667
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 9977.799999999999
$ws.Range("I132").Value = 949.5
$ws.Range("K132").Value = 2848.5
$ws.Range("M132").Value = -318.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 10008.75
$ws.Range("I36").Value = 5012.3335
$ws.Range("J36").Value = 24998
$ws.Range("K36").Value = 5012.3335
$ws.Range("L36").Value = 24998
$ws.Range("M36").Value = -4478.3335
$ws.Range("N36").Value = -26066
$ws.Range("H86").Value = 4124.75
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 4124.75
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11813.857
$ws.Range("I132").Value = 9899.333000000001
$ws.Range("J132").Value = 13249.75
$ws.Range("K132").Value = 29697.999
$ws.Range("L132").Value = 39749.25
$ws.Range("M132").Value = -27167.999
$ws.Range("N132").Value = -44809.25
$ws.Range("H134").Value = 10248.75
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 11998.333
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 35994.999
$ws.Range("M134").Value = -12465
$ws.Range("N134").Value = -41064.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 251349.5
$ws.Range("I4").Value = 334799.34
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 1004398.02
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -1004286.02
$ws.Range("N4").Value = -3224
$ws.Range("H33").Value = 55.285713
$ws.Range("I33").Value = 55.285713
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 331.714278
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -48.71427800000004
$ws.Range("N33").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 28000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 28000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 28000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -29166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2371
$ws.Range("I55").Value = 533
$ws.Range("K55").Value = 533
$ws.Range("M55").Value = -360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 15000
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224
$ws.Range("H4").Value = 19000
$ws.Range("J4").Value = 19000
$ws.Range("L4").Value = 19000
$ws.Range("N4").Value = -19226
$ws.Range("H10").Value = 663.6667
$ws.Range("J10").Value = 745
$ws.Range("L10").Value = 745
$ws.Range("N10").Value = -1083
$ws.Range("H63").Value = 26079.666
$ws.Range("I63").Value = 9995
$ws.Range("J63").Value = 34122
$ws.Range("K63").Value = 9995
$ws.Range("L63").Value = 34122
$ws.Range("M63").Value = -9371
$ws.Range("N63").Value = -35370
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 26079.666
$ws.Range("I66").Value = 9995
$ws.Range("J66").Value = 34122
$ws.Range("K66").Value = 29985
$ws.Range("L66").Value = 102366
$ws.Range("M66").Value = -26865
$ws.Range("N66").Value = -108606
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 1405
$ws.Range("I122").Value = 1441.8182
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4325.4546
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1875.4546
$ws.Range("N122").Value = -7900
